$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2058047493403694
$ws.Range("C2").Value = 0.5329815303430079
$ws.Range("J2").Value = 0.0079155672823219
$ws.Range("P2").Value = 0.1609498680738786
$ws.Range("S2").Value = 0.09234828496042216
$ws.Range("B3").Value = 0.01408450704225352
$ws.Range("C3").Value = 0.03286384976525822
$ws.Range("J3").Value = 0.02347417840375587
$ws.Range("P3").Value = 0.755868544600939
$ws.Range("S3").Value = 0.1737089201877934
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("O4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.7272727272727273
$ws.Range("S4").Value = 0.2272727272727273
$ws.Range("B6").Value = 0.08849557522123894
$ws.Range("D6").Value = 0.01327433628318584
$ws.Range("F6").Value = 0.05752212389380531
$ws.Range("J6").Value = 0.2256637168141593
$ws.Range("O6").Value = 0.02654867256637168
$ws.Range("Q6").Value = 0.1769911504424779
$ws.Range("R6").Value = 0.06194690265486726
$ws.Range("S6").Value = 0.3495575221238938
$ws.Range("B7").Value = 0.1020408163265306
$ws.Range("D7").Value = 0.03061224489795918
$ws.Range("F7").Value = 0.03061224489795918
$ws.Range("J7").Value = 0.1377551020408163
$ws.Range("O7").Value = 0.01020408163265306
$ws.Range("Q7").Value = 0.1377551020408163
$ws.Range("R7").Value = 0.07142857142857142
$ws.Range("S7").Value = 0.4795918367346939
$ws.Range("B8").Value = 0.1365853658536585
$ws.Range("D8").Value = 0.02682926829268293
$ws.Range("F8").Value = 0.03902439024390244
$ws.Range("J8").Value = 0.1219512195121951
$ws.Range("O8").Value = 0.00975609756097561
$ws.Range("Q8").Value = 0.148780487804878
$ws.Range("R8").Value = 0.0951219512195122
$ws.Range("S8").Value = 0.4219512195121951
$ws.Range("B9").Value = 0.1777777777777778
$ws.Range("D9").Value = 0.01481481481481482
$ws.Range("F9").Value = 0.01481481481481482
$ws.Range("J9").Value = 0.1037037037037037
$ws.Range("O9").Value = 0.01481481481481482
$ws.Range("Q9").Value = 0.1259259259259259
$ws.Range("R9").Value = 0.1333333333333333
$ws.Range("S9").Value = 0.4148148148148148
$ws.Range("B10").Value = 0.139005897219882
$ws.Range("D10").Value = 0.01853411962931761
$ws.Range("E10").Value = 0.0008424599831508003
$ws.Range("F10").Value = 0.08424599831508003
$ws.Range("J10").Value = 0.1086773378264532
$ws.Range("O10").Value = 0.02948609941027801
$ws.Range("Q10").Value = 0.1929233361415333
$ws.Range("R10").Value = 0.07413647851727043
$ws.Range("S10").Value = 0.3521482729570345
$ws.Range("G11").Value = 0.1466275659824047
$ws.Range("J11").Value = 0.1202346041055719
$ws.Range("K11").Value = 0.187683284457478
$ws.Range("L11").Value = 0.5249266862170088
$ws.Range("S11").Value = 0.02052785923753666
$ws.Range("G12").Value = 0.7111111111111111
$ws.Range("J12").Value = 0.1722222222222222
$ws.Range("K12").Value = 0.02222222222222222
$ws.Range("L12").Value = 0.01666666666666667
$ws.Range("S12").Value = 0.07777777777777778
$ws.Range("G13").Value = 0.8275862068965517
$ws.Range("J13").Value = 0.103448275862069
$ws.Range("S13").Value = 0.06896551724137931
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.009569377990430622
$ws.Range("H15").Value = 0.2105263157894737
$ws.Range("I15").Value = 0.06698564593301436
$ws.Range("J15").Value = 0.3110047846889952
$ws.Range("K15").Value = 0.07177033492822966
$ws.Range("M15").Value = 0.004784688995215311
$ws.Range("O15").Value = 0.04784688995215311
$ws.Range("S15").Value = 0.277511961722488
$ws.Range("F16").Value = 0.02
$ws.Range("H16").Value = 0.168
$ws.Range("I16").Value = 0.076
$ws.Range("J16").Value = 0.412
$ws.Range("K16").Value = 0.108
$ws.Range("M16").Value = 0.008
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.168
$ws.Range("F17").Value = 0.03021978021978022
$ws.Range("H17").Value = 0.1291208791208791
$ws.Range("I17").Value = 0.04945054945054945
$ws.Range("J17").Value = 0.4395604395604396
$ws.Range("K17").Value = 0.1181318681318681
$ws.Range("M17").Value = 0.01648351648351648
$ws.Range("O17").Value = 0.0576923076923077
$ws.Range("S17").Value = 0.1593406593406593
$ws.Range("F18").Value = 0.01149425287356322
$ws.Range("H18").Value = 0.1551724137931035
$ws.Range("I18").Value = 0.07471264367816093
$ws.Range("J18").Value = 0.4252873563218391
$ws.Range("K18").Value = 0.132183908045977
$ws.Range("M18").Value = 0.005747126436781609
$ws.Range("O18").Value = 0.05747126436781609
$ws.Range("S18").Value = 0.1379310344827586
$ws.Range("F19").Value = 0.03669008587041374
$ws.Range("H19").Value = 0.2076502732240437
$ws.Range("I19").Value = 0.05776736924277908
$ws.Range("J19").Value = 0.345823575331772
$ws.Range("K19").Value = 0.1178766588602654
$ws.Range("M19").Value = 0.01795472287275566
$ws.Range("N19").Value = 0.00078064012490242
$ws.Range("O19").Value = 0.0663544106167057
$ws.Range("S19").Value = 0.1491022638563622
